# Updated symbol list on Sat Feb 11 19:25:59 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "308.88"
$ws.Range("E2").Value = "0.54%"

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "40.82"
$ws.Range("E3").Value = "-0.46%"

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.116"
$ws.Range("E4").Value = "1.08%"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07618"
$ws.Range("E5").Value = "0.08%"

# Row 6
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.284"
$ws.Range("E6").Value = "0.87%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.605"
$ws.Range("E7").Value = "0.56%"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "2.464"
$ws.Range("E8").Value = "0.90%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9077"
$ws.Range("E9").Value = "0.22%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1283"
$ws.Range("E10").Value = "27.21%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1803"
$ws.Range("E11").Value = "1.89%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09037"
$ws.Range("E12").Value = "-0.34%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04300"
$ws.Range("E13").Value = "-1.97%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1044"
$ws.Range("E14").Value = "-0.90%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001257"
$ws.Range("E15").Value = "0.39%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005676"
$ws.Range("E16").Value = "-3.20%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.345"
$ws.Range("E17").Value = "-0.65%"

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.52%"

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "6.977"
$ws.Range("E19").Value = "2.57%"

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1394"
$ws.Range("E20").Value = "2.71%"

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2707"
$ws.Range("E21").Value = "-4.77%"

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04048"
$ws.Range("E22").Value = "-2.60%"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001270"
$ws.Range("E23").Value = "5.04%"

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004040"
$ws.Range("E24").Value = "-0.59%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001272"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "24.65%"

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02415"
$ws.Range("E38").Value = "0.44%"

# Row 39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05216"
$ws.Range("E39").Value = "1.60%"

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007832"
$ws.Range("E40").Value = "0.08%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006803"

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001933"
$ws.Range("E43").Value = "-0.82%"

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007340"
$ws.Range("E44").Value = "-12.48%"

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3350"
$ws.Range("E45").Value = "9.67%"

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006898"
$ws.Range("E46").Value = "8.29%"

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.10%"

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1189"
$ws.Range("E48").Value = "2,082.83%"

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.10%"

# Row 51
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.10%"
